$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 entirely (the "H 72" record), shifting all rows below up by one.
$ws.Rows.Item(2).Delete()
